# ajustes na tabela de dados
# Adds a "%" column (D) with percentage formulas (C/B) for each programmer row,
# formats it as Percent, centers the header row, and updates the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for new column D ---
$ws.Range("D1").Value = "%"

# --- Fill in Linhas (B) / Linhas Comentadas (C) with 1 for rows that only had a name ---
$ws.Range("B2:B4").Value = 1
$ws.Range("C2:C4").Value = 1
$ws.Range("B6:B11").Value = 1
$ws.Range("C6:C11").Value = 1

# --- Percentage formulas in column D for every data row ---
$ws.Range("D2:D4").Formula = "=(C2/B2)"
$ws.Range("D5").Formula = "=(C5/B5)"
$ws.Range("D6:D11").Formula = "=(C6/B6)"

# --- Number formatting: percentage style for D2:D11 ---
$ws.Range("D2:D11").Style = "Percent"

# --- Center-align the header row (A1:D1) ---
$ws.Range("A1:D1").HorizontalAlignment = -4108

# --- Update selection to H8 ---
$null = $ws.Range("H8").Select()
